# Actualización automática hashcode dom oct  6 02:05:11 CEST 2019
# Update hashcode values (column B) for the rows identified by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    "B94"  = "44213aeeab26b84a909d27da8747f1dd"
    "B95"  = "62fb3a25e5eb73fa548e78df049eeae4"
    "B98"  = "7e28e709da59e3fc566edfc13a487028"
    "B109" = "4eadddab98df18409f53e51a7d916afb"
    "B115" = "78fb34603fc974bb8815be6ff28d67f3"
    "B159" = "7efd4d5ecec095ae0b2a2e3bc16c6c20"
    "B169" = "413a0d05a619a60b898eb259c960afc0"
    "B227" = "2d01a5278488f10b9f5dd5e43c9859b6"
    "B232" = "c7017acfe56676dd01830aabf3c16619"
    "B302" = "128c4596fca9a98de68b10dcf6d5b902"
    "B339" = "4dd4c1f8cdc1fd5cc6e0107860789455"
    "B420" = "bf3569543f5afe0bd329968445d710df"
    "B464" = "f3ca3a5e106381f567089cfeb1ff5eaa"
    "B483" = "7db025c699f5ae5fc290487270fbbc2d"
    "B506" = "32b0e69ac96ccda0211b74f7e415d067"
    "B524" = "e0be8f01f61a7e46740ea82661e2c46f"
    "B558" = "48430e2174399aad2d97e1908c082c03"
    "B600" = "98a7a4c7e45a4c7f13b04e8c8f695464"
    "B624" = "23a05fa1b6ac27eb97b8412b67c6f222"
    "B626" = "cdeec3a4e361cc7e3e633c7a2be1280d"
    "B635" = "31d4b27f68ee3e27be775bef84187400"
    "B708" = "c162b077d372826d0847e23a22cd1573"
    "B827" = "4cd16c911c9d83985478f327f616afa4"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
